$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Locator Type"
$ws.Range("D2").Value = "Xpath"
$ws.Range("D3").Value = "CSS"
$ws.Range("D4").Value = "CSS"
$ws.Range("D5").Value = "CSS"
$ws.Range("D6").Value = "CSS"
$ws.Range("D7").Value = "CSS"
$ws.Range("D8").Value = "CSS"
$ws.Range("D9").Value = "CSS"
$ws.Range("D10").Value = "Xpath"
$ws.Range("D11").Value = "Xpath"
$ws.Range("D12").Value = "CSS"
$ws.Range("D13").Value = "Xpath"
$ws.Range("D14").Value = "CSS"

$ws.Range("D14").Select() | Out-Null
$ws.Columns.Item(4).ColumnWidth = 20.3
